$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.049.51"
$ws.Range("E2").Value = "  -1.48%  "
$ws.Range("D3").Value = "1.666.35"
$ws.Range("E3").Value = "  -1.26%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.88"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.61%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5098"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.86%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.004"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.21%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2669"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06392"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.65%  "
$ws.Range("E10").Value = "  -0.52%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07459"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.31%  "
$ws.Range("D12").Value = "1.695.47"
$ws.Range("E12").Value = "  +0.31%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.515"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.71%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5808"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.53%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000008513"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.82%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.08"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.92%  "
$ws.Range("D17").Value = "26.082.77"
$ws.Range("E17").Value = "  -1.61%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.924"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.47%  "
$ws.Range("E19").Value = "  -0.16%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.79"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.14%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "189.93"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.27%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.186"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.34%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.005"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.19%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "145.27"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.53%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "7.605"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.31%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1209"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.39%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.63"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.46%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06644"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +15.76%  "
$ws.Range("E29").Value = "  -1.20%  "
$ws.Range("E30").Value = "  -1.83%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.549"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.58%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.514"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.19%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.661"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.29%  "
$ws.Range("E34").Value = "  -0.25%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6133"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.372"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.46%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.689"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.38%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.383"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +8.21%  "
$ws.Range("D39").Value = "1.093.62"
$ws.Range("E39").Value = "  -0.49%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01592"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.10%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8680"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("E42").Value = "  +0.42%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "101.52"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.67%  "
$ws.Range("D44").Value = "1.813.16"
$ws.Range("E44").Value = "  -1.82%  "
$ws.Range("E45").Value = "  -1.80%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "56.31"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.16%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.007"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.23%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.067"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.33%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05226"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.15%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4287"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.77%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.003"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.27%  "

Write-Host "Applied cryptos update"